$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text number format on the Price/Volume columns for the data rows so that
# string values such as "1.030" or "0.07417" are stored verbatim (not coerced to
# numbers, which would silently drop significant trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '27.627.88'
$ws.Range("E2").Value = '  +2.99%  '

$ws.Range("D3").Value = '1.852.53'
$ws.Range("E3").Value = '  +2.20%  '

$ws.Range("E4").Value = '  +2.77%  '

$ws.Range("D5").Value = '320.73'
$ws.Range("E5").Value = '  +3.75%  '

$ws.Range("D6").Value = '1.029'
$ws.Range("E6").Value = '  +2.53%  '

$ws.Range("D7").Value = '0.4382'
$ws.Range("E7").Value = '  +1.51%  '

$ws.Range("D8").Value = '0.3756'
$ws.Range("E8").Value = '  +1.33%  '

$ws.Range("D9").Value = '0.07417'
$ws.Range("E9").Value = '  +2.40%  '

$ws.Range("D10").Value = '0.8755'
$ws.Range("E10").Value = '  +1.11%  '

$ws.Range("D11").Value = '21.54'
$ws.Range("E11").Value = '  +3.08%  '

$ws.Range("D12").Value = '1.866.14'
$ws.Range("E12").Value = '  -4.32%  '

$ws.Range("D13").Value = '5.522'
$ws.Range("E13").Value = '  +3.03%  '

$ws.Range("D14").Value = '6.703'
$ws.Range("E14").Value = '  +0.96%  '

$ws.Range("D15").Value = '0.07196'
$ws.Range("E15").Value = '  +3.98%  '

$ws.Range("D16").Value = '82.72'
$ws.Range("E16").Value = '  +2.59%  '

$ws.Range("D17").Value = '1.034'
$ws.Range("E17").Value = '  +2.59%  '

$ws.Range("D18").Value = '0.000009051'
$ws.Range("E18").Value = '  +1.48%  '

$ws.Range("D19").Value = '1.028'
$ws.Range("E19").Value = '  +2.44%  '

$ws.Range("D20").Value = '15.47'
$ws.Range("E20").Value = '  +1.66%  '

$ws.Range("D21").Value = '27.646.91'
$ws.Range("E21").Value = '  +2.92%  '

$ws.Range("D22").Value = '5.266'
$ws.Range("E22").Value = '  +1.10%  '

$ws.Range("D23").Value = '11.25'
$ws.Range("E23").Value = '  +0.48%  '

$ws.Range("D24").Value = '2.074.23'
$ws.Range("E24").Value = '  -4.23%  '

$ws.Range("D25").Value = '157.74'
$ws.Range("E25").Value = '  +2.42%  '

$ws.Range("D26").Value = '1.947'
$ws.Range("E26").Value = '  +4.14%  '

$ws.Range("D27").Value = '18.77'
$ws.Range("E27").Value = '  +2.53%  '

$ws.Range("D28").Value = '5.324'
$ws.Range("E28").Value = '  +1.79%  '

$ws.Range("D29").Value = '1.941'
$ws.Range("E29").Value = '  +2.41%  '

$ws.Range("D30").Value = '116.32'
$ws.Range("E30").Value = '  +0.90%  '

$ws.Range("D31").Value = '0.09063'
$ws.Range("E31").Value = '  +1.55%  '

$ws.Range("D32").Value = '1.212'
$ws.Range("E32").Value = '  +3.48%  '

$ws.Range("D33").Value = '0.7703'
$ws.Range("E33").Value = '  +1.61%  '

$ws.Range("D34").Value = '4.533'
$ws.Range("E34").Value = '  +2.15%  '

$ws.Range("D35").Value = '2.882'
$ws.Range("E35").Value = '  +2.57%  '

$ws.Range("D36").Value = '1.030'
$ws.Range("E36").Value = '  +2.37%  '

$ws.Range("D37").Value = '1.155'
$ws.Range("E37").Value = '  +2.04%  '

$ws.Range("D38").Value = '0.01981'
$ws.Range("E38").Value = '  +2.83%  '

$ws.Range("D39").Value = '0.05298'
$ws.Range("E39").Value = '  +1.10%  '

$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D40").Value = '0.5186'
$ws.Range("E40").Value = '  +2.16%  '

$ws.Range("B41").Value = 'MXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D41").Value = '2.819'
$ws.Range("E41").Value = '  +4.83%  '

$ws.Range("D42").Value = '0.1676'
$ws.Range("E42").Value = '  +1.56%  '

$ws.Range("D43").Value = '6.755'
$ws.Range("E43").Value = '  +2.94%  '

$ws.Range("D44").Value = '8.606'
$ws.Range("E44").Value = '  +3.63%  '

$ws.Range("B45").Value = 'Quant'
$ws.Range("C45").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D45").Value = '109.08'
$ws.Range("E45").Value = '  +2.10%  '

$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '10.63'
$ws.Range("E46").Value = '  +2.15%  '

$ws.Range("D47").Value = '1.719'
$ws.Range("E47").Value = '  +3.93%  '

$ws.Range("D48").Value = '0.4671'
$ws.Range("E48").Value = '  +2.45%  '

$ws.Range("D49").Value = '0.06391'
$ws.Range("E49").Value = '  +1.71%  '

$ws.Range("D50").Value = '1.887'
$ws.Range("E50").Value = '  +4.11%  '

$ws.Range("D51").Value = '39.65'
$ws.Range("E51").Value = '  +5.45%  '

# Restore the default cell style (clears the temporary text-number-format tag
# added above) while keeping the values/types that were just assigned.
$ws.Range("D2:E51").Style = "Normal"
